$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4.84879116695716
$ws.Range("D2").Value = 4.869771222072401
$ws.Range("E2").Value = 10.7560945941401
$ws.Range("F2").Value = 50.50267402812657
$ws.Range("G2").Value = 3.755582436368115
$ws.Range("I2").Value = 24.50450999475602
$ws.Range("J2").Value = 9.935957353701047
$ws.Range("K2").Value = 20.87932439583977
$ws.Range("M2").Value = 20.7080905607765
$ws.Range("N2").Value = 22.3180102131972
$ws.Range("B3").Value = 4.72436958968858
$ws.Range("D3").Value = 4.870079194022711
$ws.Range("E3").Value = 10.77109872029198
$ws.Range("F3").Value = 50.42106370104749
$ws.Range("G3").Value = 3.759465849973319
$ws.Range("I3").Value = 24.56268877351807
$ws.Range("J3").Value = 9.956761495249888
$ws.Range("K3").Value = 20.66533352164155
$ws.Range("M3").Value = 20.64305506517399
$ws.Range("N3").Value = 22.37739282630111
$ws.Range("B4").Value = 4.64514183132652
$ws.Range("D4").Value = 4.870548813941742
$ws.Range("E4").Value = 10.78166359961805
$ws.Range("F4").Value = 50.38286304294289
$ws.Range("G4").Value = 3.761973689669821
$ws.Range("I4").Value = 24.60176976172647
$ws.Range("J4").Value = 9.970568671844019
$ws.Range("K4").Value = 20.53872495515368
$ws.Range("M4").Value = 20.60771804480938
$ws.Range("N4").Value = 22.41579900825034
$ws.Range("B5").Value = 4.612167921280276
$ws.Range("D5").Value = 4.870810981644129
$ws.Range("E5").Value = 10.78630917428872
$ws.Range("F5").Value = 50.37029415257454
$ws.Range("G5").Value = 3.76302680559138
$ws.Range("I5").Value = 24.61853927011352
$ws.Range("J5").Value = 9.976455349916051
$ws.Range("K5").Value = 20.48838792641097
$ws.Range("M5").Value = 20.59448241758515
$ws.Range("N5").Value = 22.43193963740609
$ws.Range("B6").Value = 4.606651825497841
$ws.Range("D6").Value = 4.870858797586616
$ws.Range("E6").Value = 10.78710112874386
$ws.Range("F6").Value = 50.36838823944626
$ws.Range("G6").Value = 3.763203559616573
$ws.Range("I6").Value = 24.621374762702
$ws.Range("J6").Value = 9.977448548997074
$ws.Range("K6").Value = 20.48010695417002
$ws.Range("M6").Value = 20.59235522341961
$ws.Range("N6").Value = 22.43464938024405
$ws.Range("B7").Value = 4.644699885668883
$ws.Range("D7").Value = 4.870552062630455
$ws.Range("E7").Value = 10.78172487327004
$ws.Range("F7").Value = 50.38268139218347
$ws.Range("G7").Value = 3.761987766075217
$ws.Range("I7").Value = 24.60199250646341
$ws.Range("J7").Value = 9.97064700794988
$ws.Range("K7").Value = 20.53804093320601
$ws.Range("M7").Value = 20.60753481859089
$ws.Range("N7").Value = 22.41601470236369
$ws.Range("B8").Value = 4.806492432520637
$ws.Range("D8").Value = 4.86981933954941
$ws.Range("E8").Value = 10.76098751147587
$ws.Range("F8").Value = 50.47206617384447
$ws.Range("G8").Value = 3.756895895776539
$ws.Range("I8").Value = 24.52387231106083
$ws.Range("J8").Value = 9.942916328522463
$ws.Range("K8").Value = 20.80458110838038
$ws.Range("M8").Value = 20.68471841958873
$ws.Range("N8").Value = 22.33808207584445
$ws.Range("B9").Value = 5.100326324830294
$ws.Range("D9").Value = 4.87059761594973
$ws.Range("E9").Value = 10.73103978072681
$ws.Range("F9").Value = 50.74156804437659
$ws.Range("G9").Value = 3.747884516510303
$ws.Range("I9").Value = 24.39737360897687
$ws.Range("J9").Value = 9.896722631309016
$ws.Range("K9").Value = 21.36270448535568
$ws.Range("M9").Value = 20.87208242000512
$ws.Range("N9").Value = 22.20065642922668
$ws.Range("B10").Value = 5.300756336229026
$ws.Range("D10").Value = 4.872506156197775
$ws.Range("E10").Value = 10.71555520383344
$ws.Range("F10").Value = 50.99649486640959
$ws.Range("G10").Value = 3.741849896591527
$ws.Range("I10").Value = 24.32076741194567
$ws.Range("J10").Value = 9.867756957249
$ws.Range("K10").Value = 21.79086226087325
$ws.Range("M10").Value = 21.03100581386328
$ws.Range("N10").Value = 22.10903136605872
$ws.Range("B11").Value = 5.388373330492689
$ws.Range("D11").Value = 4.873661778517445
$ws.Range("E11").Value = 10.7099222297344
$ws.Range("F11").Value = 51.1246650524042
$ws.Range("G11").Value = 3.739230222668432
$ws.Range("I11").Value = 24.28947608894972
$ws.Range("J11").Value = 9.855655854859686
$ws.Range("K11").Value = 21.98880176362842
$ws.Range("M11").Value = 21.10775195762555
$ws.Range("N11").Value = 22.06936755915502
$ws.Range("B12").Value = 5.421024460578101
$ws.Range("D12").Value = 4.87414045604454
$ws.Range("E12").Value = 10.7079916735871
$ws.Range("F12").Value = 51.17493633300213
$ws.Range("G12").Value = 3.738256143041926
$ws.Range("I12").Value = 24.27813942443199
$ws.Range("J12").Value = 9.851227862607068
$ws.Range("K12").Value = 22.06414633262169
$ws.Range("M12").Value = 21.13743764187071
$ws.Range("N12").Value = 22.05463729919838
$ws.Range("B13").Value = 5.414016149732799
$ws.Range("D13").Value = 4.874035541932161
$ws.Range("E13").Value = 10.70839845224389
$ws.Range("F13").Value = 51.16403262827725
$ws.Range("G13").Value = 3.738465132698164
$ws.Range("I13").Value = 24.28055815815846
$ws.Range("J13").Value = 9.852174645122432
$ws.Range("K13").Value = 22.04790328000749
$ws.Range("M13").Value = 21.13101683765247
$ws.Range("N13").Value = 22.05779685295057
$ws.Range("B14").Value = 5.391070215023898
$ws.Range("D14").Value = 4.873700338035787
$ws.Range("E14").Value = 10.70975934501932
$ws.Range("F14").Value = 51.12876622478165
$ws.Range("G14").Value = 3.739149725774856
$ws.Range("I14").Value = 24.28853313331009
$ws.Range("J14").Value = 9.855288467504897
$ws.Range("K14").Value = 21.99499298736409
$ws.Range("M14").Value = 21.11018181715867
$ws.Range("N14").Value = 22.06814989110305
$ws.Range("B15").Value = 5.376946030740623
$ws.Range("D15").Value = 4.87350035650529
$ws.Range("E15").Value = 10.71061929460266
$ws.Range("F15").Value = 51.10739001598547
$ws.Range("G15").Value = 3.739571390982808
$ws.Range("I15").Value = 24.29348483595225
$ws.Range("J15").Value = 9.857215878609047
$ws.Range("K15").Value = 21.96263261631481
$ws.Range("M15").Value = 21.09750046824653
$ws.Range("N15").Value = 22.07452912219715
$ws.Range("B16").Value = 5.294957356162009
$ws.Range("D16").Value = 4.872436394780564
$ws.Range("E16").Value = 10.71595169423363
$ws.Range("F16").Value = 50.98836265012557
$ws.Range("G16").Value = 3.742023614217506
$ws.Range("I16").Value = 24.32288405316488
$ws.Range("J16").Value = 9.868569415659838
$ws.Range("K16").Value = 21.7779845265275
$ws.Range("M16").Value = 21.02607845540059
$ws.Range("N16").Value = 22.11166403393847
$ws.Range("B17").Value = 5.243736605919209
$ws.Range("D17").Value = 4.871857117333898
$ws.Range("E17").Value = 10.71958409245096
$ws.Range("F17").Value = 50.91845610395605
$ws.Range("G17").Value = 3.743560037725295
$ws.Range("I17").Value = 24.34183143909304
$ws.Range("J17").Value = 9.875809749093785
$ws.Range("K17").Value = 21.66547316150163
$ws.Range("M17").Value = 20.98339200370864
$ws.Range("N17").Value = 22.13496133619344
$ws.Range("B18").Value = 5.213941567389979
$ws.Range("D18").Value = 4.871551001127292
$ws.Range("E18").Value = 10.72180619013921
$ws.Range("F18").Value = 50.87939756653274
$ws.Range("G18").Value = 3.744455567765919
$ws.Range("I18").Value = 24.35306422892214
$ws.Range("J18").Value = 9.880075444620323
$ws.Range("K18").Value = 21.60106181669909
$ws.Range("M18").Value = 20.9592598981354
$ws.Range("N18").Value = 22.14855119409324
$ws.Range("B19").Value = 5.203796580363479
$ws.Range("D19").Value = 4.871452012509448
$ws.Range("E19").Value = 10.72258138023331
$ws.Range("F19").Value = 50.86637103506997
$ws.Range("G19").Value = 3.744760812386637
$ws.Range("I19").Value = 24.35692491735766
$ws.Range("J19").Value = 9.881537132087004
$ws.Range("K19").Value = 21.57930717128852
$ws.Range("M19").Value = 20.95116180434761
$ws.Range("N19").Value = 22.15318511061798
$ws.Range("B20").Value = 5.249223850944506
$ws.Range("D20").Value = 4.871915982776406
$ws.Range("E20").Value = 10.71918367181913
$ws.Range("F20").Value = 50.92577889673073
$ws.Range("G20").Value = 3.743395260330726
$ws.Range("I20").Value = 24.33977980325906
$ws.Range("J20").Value = 9.875028526190713
$ws.Range("K20").Value = 21.67741939132924
$ws.Range("M20").Value = 20.98789269576342
$ws.Range("N20").Value = 22.13246165165361
$ws.Range("B21").Value = 5.39782442602289
$ws.Range("D21").Value = 4.873797683053391
$ws.Range("E21").Value = 10.7093541244233
$ws.Range("F21").Value = 51.1390778690757
$ws.Range("G21").Value = 3.738948158370754
$ws.Range("I21").Value = 24.28617676667265
$ws.Range("J21").Value = 9.854369672931744
$ws.Range("K21").Value = 22.01052398996623
$ws.Range("M21").Value = 21.11628478016163
$ws.Range("N21").Value = 22.06510109775055
$ws.Range("B22").Value = 5.4918625140296
$ws.Range("D22").Value = 4.875266739036654
$ws.Range("E22").Value = 10.70411024063467
$ws.Range("F22").Value = 51.28858945352989
$ws.Range("G22").Value = 3.736146201103714
$ws.Range("I22").Value = 24.25413289927522
$ws.Range("J22").Value = 9.841767965787154
$ws.Range("K22").Value = 22.23046482481375
$ws.Range("M22").Value = 21.20382281588649
$ws.Range("N22").Value = 22.02276467865784
$ws.Range("B23").Value = 5.44195929325454
$ws.Range("D23").Value = 4.874460871054968
$ws.Range("E23").Value = 10.70680113669119
$ws.Range("F23").Value = 51.20787421803526
$ws.Range("G23").Value = 3.737632135778234
$ws.Range("I23").Value = 24.27096146062878
$ws.Range("J23").Value = 9.848411453988646
$ws.Range("K23").Value = 22.11289542388364
$ws.Range("M23").Value = 21.15677600689263
$ws.Range("N23").Value = 22.04520615372806
$ws.Range("B24").Value = 5.246744148795171
$ws.Range("D24").Value = 4.871889285846745
$ws.Range("E24").Value = 10.71936428525288
$ws.Range("F24").Value = 50.92246473583783
$ws.Range("G24").Value = 3.74346971812085
$ws.Range("I24").Value = 24.34070628986362
$ws.Range("J24").Value = 9.875381395816781
$ws.Range("K24").Value = 21.67201764564613
$ws.Range("M24").Value = 20.98585665800363
$ws.Range("N24").Value = 22.13359114872531
$ws.Range("B25").Value = 5.023472824942586
$ws.Range("D25").Value = 4.870151331701203
$ws.Range("E25").Value = 10.73799551387197
$ws.Range("F25").Value = 50.65861660903443
$ws.Range("G25").Value = 3.75021887181987
$ws.Range("I25").Value = 24.42873089487084
$ws.Range("J25").Value = 9.90834470523887
$ws.Range("K25").Value = 21.20827860428775
$ws.Range("M25").Value = 20.81760647947141
$ws.Range("N25").Value = 22.23618994728552
